$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty "Zeitdauer Robin" (column E) durations for the
# last work packages (Installationsanleitung, Deployment/Testen,
# Dokumentation) that were left blank.
$ws.Range("E34").Value = 5
$ws.Range("E35").Value = 0.5
$ws.Range("E36").Value = 4.5

# Move the current selection to reflect where the author ended up working.
$ws.Range("M23:M24").Select()
